# Apply the "dhydro_support_hydrolib-core" content update (deployed to
# `develop` with MkDocs 1.3.0 and mike 1.1.2):
#   - Observation station file (new) is now supported as of 0.3.0,
#     backed by hydrolib.core.io.obs.models.ObservationPointModel
#   - RainfallRunoffModel moved from hydrolib.core.io.fnm.models to
#     hydrolib.core.io.rr.models
#   - BuiModel moved from hydrolib.core.io.bui.models to
#     hydrolib.core.io.rr.meteo.models

$wb = $excel.ActiveWorkbook
$wsSource = $wb.Worksheets.Item("Source table")

# Row 48: Observation station file (new)
$wsSource.Range("B48").Value = "X"
$wsSource.Range("C48").Value = "X"
$wsSource.Range("D48").Value = "0.3.0"
$wsSource.Range("E48").Value = "hydrolib.core.io.obs.models"
$wsSource.Range("F48").Value = "ObservationPointModel"

# Row 59: Main sobek_3b.fnm / RainfallRunoffModel
$wsSource.Range("E59").Value = "hydrolib.core.io.rr.models"
$wsSource.Range("G59").Value = "Used to be in hydrolib.core.io.fnm.models before 0.3.0"

# Row 60: Rainfall .bui file / BuiModel
$wsSource.Range("E60").Value = "hydrolib.core.io.rr.meteo.models"
$wsSource.Range("G60").Value = "Used to be in hydrolib.core.io.bui.models before 0.3.0"

# Reflect the author's on-screen navigation: they scrolled around the
# "Source table" sheet (ending with B29 selected there) ...
[void]$wsSource.Activate()
[void]$wsSource.Range("B29").Select()

# ... and the "FM mkdocs table" sheet (ending with A53 selected there) ...
$wsFM = $wb.Worksheets.Item("FM mkdocs table")
[void]$wsFM.Activate()
[void]$wsFM.Range("A53").Select()

# ... leaving "FM mkdocs table" as the active tab when the file was saved.
